# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be written as TEXT (matching the
# workbook's existing inline-string cells), while leaving the cell's style
# index untouched (Excel would otherwise stamp a "Text" number format on it).
function Set-TextValue($a1, $value) {
    $r = $ws.Range($a1)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "68.573.63"
$ws.Range("E2").Value = "  -1.02%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.855.16"
$ws.Range("E3").Value = "  -2.20%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "521.69"
$ws.Range("E5").Value = "  +6.01%  "

# Row 6 - Solana
Set-TextValue "D6" "140.95"
$ws.Range("E6").Value = "  -4.45%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -2.22%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.20%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -2.81%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -5.15%  "

# Row 11 - ShibaInu
$ws.Range("E11").Value = "  -8.11%  "

# Row 12 - Avalanche
Set-TextValue "D12" "41.59"
$ws.Range("E12").Value = "  -3.75%  "

# Row 13 - Polkadot
Set-TextValue "D13" "10.33"
$ws.Range("E13").Value = "  -0.98%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "4.474.65"
$ws.Range("E14").Value = "  -2.16%  "

# Row 15 - Chainlink
Set-TextValue "D15" "21.39"
$ws.Range("E15").Value = "  +7.70%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.870.08"
$ws.Range("E16").Value = "  -1.81%  "

# Row 17 - Uniswap
Set-TextValue "D17" "14.09"
$ws.Range("E17").Value = "  -1.55%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -2.12%  "

# Row 19 - Polygon
$ws.Range("E19").Value = "  +1.58%  "

# Row 20 - WrappedBTC
Set-TextValue "D20" "68.591.54"
$ws.Range("E20").Value = "  -1.20%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "414.66"

# Row 22 - ImmutableX
Set-TextValue "D22" "3.48"
$ws.Range("E22").Value = "  +0.76%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("E23").Value = "  -2.89%  "

# Row 24 - RenderToken
Set-TextValue "D24" "12.00"
$ws.Range("E24").Value = "  -0.54%  "

# Row 25 - Litecoin
Set-TextValue "D25" "86.56"
$ws.Range("E25").Value = "  -2.54%  "

# Row 26 - PancakeSwap
Set-TextValue "D26" "3.99"
$ws.Range("E26").Value = "  +5.24%  "

# Row 27 - Filecoin
Set-TextValue "D27" "10.46"
$ws.Range("E27").Value = "  -5.81%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "35.42"
$ws.Range("E28").Value = "  -4.57%  "

# Row 29 - Cosmos
Set-TextValue "D29" "13.32"
$ws.Range("E29").Value = "  -0.29%  "

# Row 30 - Bittensor
Set-TextValue "D30" "677.42"
$ws.Range("E30").Value = "  -4.07%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "6.94"
$ws.Range("E31").Value = "  +14.22%  "

# Row 32 - Hedera
Set-TextValue "D32" "0.124"

# Row 33 - Toncoin
Set-TextValue "D33" "2.79"
$ws.Range("E33").Value = "  -3.51%  "

# Row 34 - OKB
Set-TextValue "D34" "66.80"
$ws.Range("E34").Value = "  +8.55%  "

# Row 35 - TheGraph
Set-TextValue "D35" "0.450"
$ws.Range("E35").Value = "  -2.95%  "

# Row 36 - PEPE
Set-TextValue "D36" "0.0₃0847"
$ws.Range("E36").Value = "  -7.18%  "

# Row 37 - InjectiveProtocol
Set-TextValue "D37" "39.42"
$ws.Range("E37").Value = "  -3.59%  "

# Row 38 - ThetaToken
Set-TextValue "D38" "3.45"
$ws.Range("E38").Value = "  +12.88%  "

# Row 39 - Kaspa
Set-TextValue "D39" "0.148"
$ws.Range("E39").Value = "  -1.96%  "

# Row 40 - Dai
Set-TextValue "D40" "0.999"
$ws.Range("E40").Value = "  +0.02%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  -0.14%  "

# Row 42 - now Fetch.AI (was VeChain)
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D42" "2.89"
$ws.Range("E42").Value = "  -1.33%  "

# Row 43 - now WEMIXToken (was Fetch.AI)
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D43" "3.17"
$ws.Range("E43").Value = "  +5.47%  "

# Row 44 - now VeChain (was WEMIXToken)
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D44" "0.0474"
$ws.Range("E44").Value = "  -3.33%  "

# Row 45 - ApeXProtocol
$ws.Range("E45").Value = "  +1.80%  "

# Row 46 - Stellar
$ws.Range("E46").Value = "  -1.71%  "

# Row 47 - FLOKI
Set-TextValue "D47" "0.000281"
$ws.Range("E47").Value = "  +16.32%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  +0.31%  "

# Row 49 - LidoDAOToken
$ws.Range("E49").Value = "  -3.49%  "

# Row 50 - now Monero (was THORChain)
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D50" "143.87"
$ws.Range("E50").Value = "  -0.16%  "

# Row 51 - now THORChain (was Monero)
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D51" "8.74"
$ws.Range("E51").Value = "  +3.84%  "
